$d = $word.ActiveDocument

function Get-ParaRangeByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p.Range
        }
    }
    return $null
}

# --- Paragraph: "NO. DE PEDIDO: ... RECURSO: " ---
$para4 = Get-ParaRangeByText("NO. DE PEDIDO:")
$xmlPara4 = "<w:p w:rsidR=`"00AE776F`" w:rsidRDefault=`"00015631`"><w:pPr><w:spacing w:after=`"60`"/><w:ind w:left=`"720`" w:hanging=`"10`"/><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`" w:eastAsia=`"es-MX`"/></w:rPr><mc:AlternateContent><mc:Choice Requires=`"wps`"><w:drawing><wp:anchor distT=`"0`" distB=`"0`" distL=`"114300`" distR=`"114300`" simplePos=`"0`" relativeHeight=`"251663360`" behindDoc=`"0`" locked=`"0`" layoutInCell=`"1`" allowOverlap=`"1`"><wp:simplePos x=`"0`" y=`"0`"/><wp:positionH relativeFrom=`"margin`"><wp:posOffset>4321175</wp:posOffset></wp:positionH><wp:positionV relativeFrom=`"paragraph`"><wp:posOffset>26670</wp:posOffset></wp:positionV><wp:extent cx=`"1552575`" cy=`"295275`"/><wp:effectExtent l=`"0`" t=`"0`" r=`"9525`" b=`"9525`"/><wp:wrapNone/><wp:docPr id=`"3`" name=`"Cuadro de texto 3`"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a=`"http://schemas.openxmlformats.org/drawingml/2006/main`"><a:graphicData uri=`"http://schemas.microsoft.com/office/word/2010/wordprocessingShape`"><wps:wsp><wps:cNvSpPr txBox=`"1`"/><wps:spPr><a:xfrm><a:off x=`"0`" y=`"0`"/><a:ext cx=`"1552575`" cy=`"295275`"/></a:xfrm><a:prstGeom prst=`"rect`"><a:avLst/></a:prstGeom><a:solidFill><a:schemeClr val=`"lt1`"/></a:solidFill><a:ln w=`"6350`"><a:noFill/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`">FOLIO: </w:t></w:r><w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>`${folio}</w:t></w:r><w:bookmarkEnd w:id=`"0`"/></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot=`"0`" spcFirstLastPara=`"0`" vertOverflow=`"overflow`" horzOverflow=`"overflow`" vert=`"horz`" wrap=`"square`" lIns=`"91440`" tIns=`"45720`" rIns=`"91440`" bIns=`"45720`" numCol=`"1`" spcCol=`"0`" rtlCol=`"0`" fromWordArt=`"0`" anchor=`"t`" anchorCtr=`"0`" forceAA=`"0`" compatLnSpc=`"1`"><a:prstTxWarp prst=`"textNoShape`"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom=`"margin`"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype id=`"_x0000_t202`" coordsize=`"21600,21600`" o:spt=`"202`" path=`"m,l,21600r21600,l21600,xe`"><v:stroke joinstyle=`"miter`"/><v:path gradientshapeok=`"t`" o:connecttype=`"rect`"/></v:shapetype><v:shape id=`"Cuadro de texto 3`" o:spid=`"_x0000_s1027`" type=`"#_x0000_t202`" style=`"position:absolute;left:0;text-align:left;margin-left:340.25pt;margin-top:2.1pt;width:122.25pt;height:23.25pt;z-index:251663360;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:margin;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-width-relative:margin;v-text-anchor:top`" fillcolor=`"white [3201]`" stroked=`"f`" strokeweight=`".5pt`"><v:textbox><w:txbxContent><w:p><w:bookmarkStart w:id=`"1`" w:name=`"_GoBack`"/><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`">FOLIO: </w:t></w:r><w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t>`${folio}</w:t></w:r><w:bookmarkEnd w:id=`"1`"/></w:p></w:txbxContent></v:textbox><w10:wrap anchorx=`"margin`"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`">                                     </w:t></w:r><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`">                                            </w:t></w:r><w:r><w:rPr><w:sz w:val=`"21`"/><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>"
$para4.InsertXML($xmlPara4)

# --- Paragraph: "FOLIO: `${folio}" becomes empty ---
$para5 = Get-ParaRangeByText("FOLIO:")
$xmlPara5 = "<w:p w:rsidR=`"00AE776F`" w:rsidRDefault=`"00015631`"><w:pPr><w:spacing w:after=`"60`"/><w:ind w:left=`"720`" w:hanging=`"10`"/><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr></w:pPr></w:p>"
$para5.InsertXML($xmlPara5)

# --- Table cell paragraph "    `${dia}" (merge 3 runs into 1) ---
$paraDia = Get-ParaRangeByText('${')
$xmlParaDia = "<w:p w:rsidR=`"00AE776F`" w:rsidRDefault=`"00015631`"><w:pPr><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"es-MX`"/></w:rPr><w:t xml:space=`"preserve`">    `${dia}</w:t></w:r></w:p>"
$paraDia.InsertXML($xmlParaDia)
